$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = '[''4'', ''4:30'']'
$ws.Range("D6").Value = '[''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D7").Value = '[''12:30'', ''1'', ''1:30''] and  [''3'', ''3:30'', ''4'', ''4:30'']  '
$ws.Range("D10").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30''] and  [''3'', ''3:30'', ''4'', ''4:30'']  '
$ws.Range("D11").Value = '[''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D13").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D16").Value = '[''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D17").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D19").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D20").Value = '[''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30''] and  [''4'', ''4:30'']  '
$ws.Range("D21").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D23").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D24").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D25").Value = '[''4'', ''4:30'']'
$ws.Range("D26").Value = '[''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D27").Value = '[''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30''] and  [''4'', ''4:30'']  '
$ws.Range("D29").Value = '[''10'', ''10:30''] and  [''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']  '
$ws.Range("D30").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D31").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2''] and  [''4'', ''4:30'']  '
$ws.Range("D32").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D33").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D34").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D35").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D36").Value = '[''10'', ''10:30'', ''11'', ''11:30'', ''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D37").Value = '[''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D38").Value = '[''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D39").Value = '[''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D40").Value = '[''10'', ''10:30'', ''11'', ''11:30''] and  [''3'', ''3:30'', ''4'', ''4:30'']  '
$ws.Range("D41").Value = '[''3:30'', ''4'', ''4:30'']'
$ws.Range("D42").Value = '[''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D43").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D44").Value = '[''12'', ''12:30'', ''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
$ws.Range("D45").Value = '[''1'', ''1:30'', ''2'', ''2:30'', ''3'', ''3:30'', ''4'', ''4:30'']'
